# Test Cases: adding status column
# Adds a new "Status" column (J) to the test cases sheet, populated with
# "Automated" / "Manual" for each existing test-case row, formatted to
# match the workbook's existing header/body styling, and tidies up the
# sheet's used range (drops the two stray fully-empty trailing rows that
# Excel round-tripped down at 1048575/1048576).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (J1): same look as the other header cells -----------------
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J1").Value = "Status"

# --- Body cells (J2:J40): centered, wrapped "Automated"/"Manual" values ----
$ws.Range("A2").Copy()
$ws.Range("J2:J40").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("J2:J40").HorizontalAlignment = -4108   # xlCenter
$ws.Range("J2:J40").VerticalAlignment = -4108     # xlCenter
$ws.Range("J2:J40").WrapText = $true

$statusByRow = @{
    2  = "Automated"
    3  = "Manual"
    4  = "Manual"
    5  = "Automated"
    6  = "Manual"
    7  = "Automated"
    8  = "Manual"
    9  = "Manual"
    10 = "Automated"
    11 = "Automated"
    12 = "Automated"
    13 = "Automated"
    14 = "Manual"
    15 = "Automated"
    16 = "Manual"
    17 = "Manual"
    18 = "Automated"
    19 = "Manual"
    20 = "Manual"
    21 = "Manual"
    22 = "Automated"
    23 = "Manual"
    24 = "Automated"
    25 = "Manual"
    26 = "Manual"
    27 = "Automated"
    28 = "Automated"
    29 = "Automated"
    30 = "Manual"
    31 = "Manual"
    32 = "Manual"
    33 = "Manual"
    34 = "Manual"
    35 = "Manual"
    36 = "Manual"
    37 = "Automated"
    38 = "Manual"
    39 = "Manual"
    40 = "Manual"
}

foreach ($row in $statusByRow.Keys) {
    $ws.Cells.Item($row, 10).Value = $statusByRow[$row]
}

# Row 17 gained a stray empty-but-styled cell at I37 in the same pass.
$ws.Range("J2").Copy()
$ws.Range("I37").PasteSpecial(-4122)   # xlPasteFormats

# --- Column width for the new column ---------------------------------------
$ws.Columns.Item(10).ColumnWidth = 14.0

# --- Drop the two stray fully-empty trailing rows ---------------------------
$ws.Range("A1048575:A1048576").EntireRow.Delete()

# --- Selection / view state to match the saved workbook ---------------------
$ws.Range("J40").Select()
